# Edit: Add 2 more bills (Bill_16 and Bill_17) to the "Expenses" sheet.
# This inserts two new data rows (new row 19 and row 20) just above the
# existing "Total" summary block, which in turn gets pushed down by two
# rows (old rows 19/20/21 become rows 21/22/23), and updates the related
# formulas, shared strings, column width and selection/view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Expenses")

# ---------------------------------------------------------------------
# Step 1: relocate the existing summary rows (19,20,21) down to (21,22,23)
# Work from the bottom row upward so we never overwrite data we still need.
# ---------------------------------------------------------------------

# --- old row 21 (Received / Balance) -> new row 23 ---
$ws.Range("A23").Value = $ws.Range("A21").Value2
$ws.Range("B23").Value = $ws.Range("B21").Value2
$ws.Range("C23").Value = $ws.Range("C21").Value2
$ws.Range("D23").Value = $ws.Range("D21").Value2
$ws.Range("E23").Formula = "=Internal!B3+Internal!B6"
$ws.Range("F23").Value = $ws.Range("F21").Value2
$ws.Range("G23").Value = $ws.Range("G21").Value2
$ws.Range("H23").Formula = "=E23-E22"

# --- old row 20 (Total Petrol + Expenses) -> new row 22 ---
$ws.Range("A22").Value = $ws.Range("A20").Value2
$ws.Range("B22").Value = $ws.Range("B20").Value2
$ws.Range("C22").Value = $ws.Range("C20").Value2
$ws.Range("D22").Value = $ws.Range("D20").Value2
$ws.Range("E22").Formula = "=E21+H21"
$ws.Range("F22").Value = $ws.Range("F20").Value2
$ws.Range("G22").Value = $ws.Range("G20").Value2
$ws.Range("H22").Value = $ws.Range("H20").Value2

# --- old row 19 (Total) -> new row 21 ---
$ws.Range("A21").Value = $ws.Range("A19").Value2
$ws.Range("B21").Value = $ws.Range("B19").Value2
$ws.Range("C21").Value = $ws.Range("C19").Value2
$ws.Range("D21").Value = $ws.Range("D19").Value2
$ws.Range("E21").Formula = "=SUM(E2:E19)"
$ws.Range("F21").Value = $ws.Range("F19").Value2
$ws.Range("G21").Value = $ws.Range("G19").Value2
$ws.Range("H21").Formula = "=SUM(H2:H18)"

# Copy formatting for the relocated summary rows from their previous
# locations (rows 19,20,21) onto their new locations (21,22,23).
$ws.Range("A19:J19").Copy() | Out-Null
$ws.Range("A21:J21").PasteSpecial(-4122) | Out-Null
$ws.Range("A20:J20").Copy() | Out-Null
$ws.Range("A22:J22").PasteSpecial(-4122) | Out-Null
$ws.Range("A21:J21").Copy() | Out-Null
$ws.Range("A23:J23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Step 2: write the two new bill rows into (now-free) rows 19 and 20
# ---------------------------------------------------------------------

# Row 19 - Bill_16 (Rajiv Electronics Bill)
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = 44389
$ws.Range("C19").Value = "Rajiv Electronics Bill"
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = 1105
$ws.Range("F19").Value = "Akshay"
$ws.Range("G19").Value = ""
$ws.Range("H19").Value = ""
$ws.Range("I19").Value = ""
$ws.Range("J19").Value = ""

# Row 20 - Bill_17 (Bus Ticket Bill and Auto travel)
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = 44390
$ws.Range("C20").Value = "Bus Ticket Bill and Auto travel"
$ws.Range("D20").Value = 17
$ws.Range("E20").Value = 2205
$ws.Range("F20").Value = "Akshay"
$ws.Range("G20").Value = ""
$ws.Range("H20").Value = 500
$ws.Range("I20").Value = "(Auto500)"
$ws.Range("J20").Value = ""

# Apply the same formatting used by the row above (row 18, the last
# original bill row) to the two freshly written rows.
$ws.Range("A18:J18").Copy() | Out-Null
$ws.Range("A19:J20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 20 wraps onto two lines (long text in column C), so its height
# grows to 30 points, same as the other two-line rows in the sheet.
$ws.Rows.Item(20).RowHeight = 30

# ---------------------------------------------------------------------
# Step 3: misc sheet-level formatting updates
# ---------------------------------------------------------------------

# Give column I (new bill remark column) a bit more width.
$ws.Columns.Item(9).ColumnWidth = 11

# Update the view: scrolled position and current selection.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 12
$ws.Range("G20").Select()
